# DDAf_2023_Tableau_annexe_Tab10.xlsx update
# - Soudan du Sud and Cabo Verde are (re)marked as resource-rich countries
#   (country name gains a trailing "*"); Nigeria loses that marker.
# - The "resource-rich" row shading follows the marker: Soudan du Sud's row
#   gains the light-blue fill, Nigeria's row loses it (Cabo Verde's row was
#   already shaded).
# - All aggregate rows whose grouping depends on the resource-rich flag
#   (the various "… pays riches en ressources …" / "… exclus" / income-group
#   / fragile-state rollups) are refreshed to their recomputed averages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name (shared string) corrections -----------------------------
$ws.Range("B34").Value = "Soudan du Sud*"
$ws.Range("B48").Value = "Cabo Verde*"
$ws.Range("B57").Value = "Nigeria"

# --- Row shading: mirror the resource-rich flag ----------------------------
# Soudan du Sud (row 34) becomes resource-rich -> apply the shaded format
# already used by other resource-rich country rows (e.g. row 20).
$ws.Range("B20:J20").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122)   # xlPasteFormats

# Nigeria (row 57) is no longer resource-rich -> apply the plain (unshaded)
# format used by other non-resource-rich country rows (e.g. row 22).
$ws.Range("B22:J22").Copy()
$ws.Range("B57:J57").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Recalculated aggregate rows -------------------------------------------
$ws.Range("C69").Value = 0.50900000000000001
$ws.Range("D69").Value = 0.91128571428571004
$ws.Range("E69").Value = 0.52200000000000002
$ws.Range("F69").Value = 38.071428571428598
$ws.Range("G69").Value = 54.1
$ws.Range("H69").Value = 30.8857142857143
$ws.Range("I69").Value = 39.9428571428571
$ws.Range("J69").Value = 22.714285714285701

$ws.Range("C77").Value = 0.76790909090908999
$ws.Range("D77").Value = 0.97590909090908995
$ws.Range("E77").Value = 0.38536363636364002
$ws.Range("F77").Value = 21.9444444444444
$ws.Range("G77").Value = 23.954545454545499
$ws.Range("H77").Value = 25.490909090909099
$ws.Range("I77").Value = 24.255555555555599
$ws.Range("J77").Value = 17.1636363636364

$ws.Range("C80").Value = 0.59640000000000004
$ws.Range("D80").Value = 0.89177777777778
$ws.Range("F80").Value = 42.24
$ws.Range("G80").Value = 54.4
$ws.Range("H80").Value = 29.285714285714299
$ws.Range("I80").Value = 54.985714285714302
$ws.Range("J80").Value = 38.287500000000001

$ws.Range("C82").Value = 0.55244186046512
$ws.Range("D82").Value = 0.89815384615384997
$ws.Range("F82").Value = 40.389189189189203
$ws.Range("G82").Value = 51.986363636363699
$ws.Range("H82").Value = 29.6794871794872
$ws.Range("I82").Value = 38.7358974358974
$ws.Range("J82").Value = 33.857500000000002

$ws.Range("C84").Value = 0.47447826086957001
$ws.Range("D84").Value = 0.86833333333332996
$ws.Range("E84").Value = 0.58461111111110997
$ws.Range("F84").Value = 40.75
$ws.Range("G84").Value = 51.2916666666667
$ws.Range("H84").Value = 32.122727272727303
$ws.Range("I84").Value = 43.886363636363697
$ws.Range("J84").Value = 30.927272727272701

$ws.Range("C86").Value = 0.59531818181817997
$ws.Range("D86").Value = 0.90285714285714003
$ws.Range("E86").Value = 0.52126315789473998
$ws.Range("F86").Value = 43.068750000000001
$ws.Range("G86").Value = 55.531818181818203
$ws.Range("H86").Value = 27.1105263157895
$ws.Range("I86").Value = 43.431578947368401
$ws.Range("J86").Value = 39.594999999999999

$ws.Range("C87").Value = 0.65893548387097001
$ws.Range("D87").Value = 0.93008000000000002
$ws.Range("E87").Value = 0.43165384615385
$ws.Range("F87").Value = 36.243478260869601
$ws.Range("G87").Value = 42.907142857142901
$ws.Range("H87").Value = 33.648148148148202
$ws.Range("I87").Value = 32.880000000000003
$ws.Range("J87").Value = 31.866666666666699

$ws.Range("C89").Value = 0.76613636363635995
$ws.Range("D89").Value = 0.96439473684210997
$ws.Range("E89").Value = 0.31865789473683998
$ws.Range("F89").Value = 25.196551724137901
$ws.Range("G89").Value = 26.610810810810801
$ws.Range("H89").Value = 28.1514285714286
$ws.Range("I89").Value = 23.163333333333298
$ws.Range("J89").Value = 21.875

$ws.Range("C90").Value = 0.88503571428571004
$ws.Range("D90").Value = 0.98065999999999998
$ws.Range("E90").Value = 0.1401568627451
$ws.Range("F90").Value = 19.480487804878099
$ws.Range("G90").Value = 26.092452830188702
$ws.Range("H90").Value = 21.207317073170699
$ws.Range("I90").Value = 15.875
$ws.Range("J90").Value = 24.419565217391298

$ws.Range("C94").Value = 0.72055172413793001
$ws.Range("D94").Value = 0.96633333333333005
$ws.Range("E94").Value = 0.39200000000000002
$ws.Range("G94").Value = 32.005000000000003
$ws.Range("I94").Value = 24.887499999999999
$ws.Range("J94").Value = 26.422222222222199

$ws.Range("C97").Value = 0.51367567567567995
$ws.Range("D97").Value = 0.88193939393939003
$ws.Range("E97").Value = 0.57425000000000004
$ws.Range("F97").Value = 41.6645161290323
$ws.Range("G97").Value = 55.065789473684198
$ws.Range("H97").Value = 31.569696969696999
$ws.Range("I97").Value = 45.103030303030302
$ws.Range("J97").Value = 34.5828571428572

$ws.Range("C98").Value = 0.61365000000000003
$ws.Range("D98").Value = 0.86170588235293999
$ws.Range("G98").Value = 48.984210526315799
$ws.Range("H98").Value = 35.941176470588204
$ws.Range("I98").Value = 37.581249999999997
